$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: swap the date value and update the refcode
$ws.Range("E2").Value = "2021-05-14 17:39:41"
$ws.Range("F2").Value = "Third"

# Row 3: swap the date value and update the refcode
$ws.Range("E3").Value = "2021-05-24 10:49:41"
$ws.Range("F3").Value = "Fourth"

# Update the active selection to E3
$ws.Range("E3").Select() | Out-Null
